$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    10 = "Chikkamagaluru (Chikmagalur)"
    11 = "Chikkamagaluru (Chikmagalur)"
    15 = "Bagalkot"
    16 = "Chikkamagaluru (Chikmagalur)"
    17 = "Ballari (Bellary)"
    18 = "Ballari (Bellary)"
    21 = "Ballari (Bellary)"
    24 = "Chikkamagaluru (Chikmagalur)"
    32 = "Ballari (Bellary)"
    33 = "Chikkamagaluru (Chikmagalur)"
    34 = "Kalaburagi (Gulbarga)"
    37 = "Kalaburagi (Gulbarga)"
    38 = "Shivamogga (Shimoga)"
    40 = "Chikkamagaluru (Chikmagalur)"
    44 = "Shivamogga (Shimoga)"
    45 = "Chikkamagaluru (Chikmagalur)"
    50 = "Chikkamagaluru (Chikmagalur)"
    51 = "Shivamogga (Shimoga)"
    52 = "Shivamogga (Shimoga)"
    53 = "Shivamogga (Shimoga)"
    54 = "Ballari (Bellary)"
    57 = "Ballari (Bellary)"
    58 = "Vijayapura (Bijapur)"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
